# Update SubRES Excel templates
# SubRES_TMPL/SubRes_2.1_ELC_EIA_AEO2022.xlsx
#
# - Column B (and K, which mirrors B) process-id formula now also folds in
#   the AH ("THM"/"RNW"/"STG" group) column, e.g. E_COA_CAV_01 -> E_THM_COA_CAV_01
# - A couple of column widths widened to fit the longer process ids / labels
# - Solar PV capacity factors (X25/X26) simplified from 0.35 down to 0.11
# - Cosmetic: active selection left on X27

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ELC_NEW_GENERATORS")

# --- 1. Rebuild the process-id formula in column B -------------------------
# Row 4 is edited first as a standalone formula (mirrors how the author
# edited a single cell before filling it down through the rest of the
# table), then the fill covers B5:B27 and re-establishes the shared
# formula group the same way Excel does when you copy a formula down.
$ws.Range("B4").Formula = '="E_"&AH4&"_"&AI4&"_"&AJ4&"_"&AL4'
$ws.Range("B5:B27").Formula = '="E_"&AH5&"_"&AI5&"_"&AJ5&"_"&AL5'

# K4:K27 simply mirror column B ("=B4" etc.) and recalc automatically from
# the formulas above, exactly as they did before the edit.

# --- 2. Column width updates -------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 22.140625    # B  (process id)
$ws.Columns.Item(11).ColumnWidth = 22.140625   # K  (process id, mirrors B)
$ws.Columns.Item(27).ColumnWidth = 26.7109375  # AA
$ws.Columns.Item(28).ColumnWidth = 26.7109375  # AB (now matches AA)
$ws.Columns.Item(37).ColumnWidth = 29.7109375  # AK

# --- 3. Simplify Solar PV capacity factors ------------------------------
$ws.Range("X25").Value = 0.11
$ws.Range("X26").Value = 0.11

# --- 4. Leave selection where the author left it ------------------------
$ws.Range("X27").Select()
